$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Joel / Adam / ydg@hotmail.com / Selenium Automation / Others
$ws.Range("A3").Value = "Joel"
$ws.Range("B3").Value = "Adam"
$ws.Range("C3").Value = "ydg@hotmail.com"

# Fill the rest of column A first (matches the original author's entry order)
$ws.Range("A4").Value = "Neel"
$ws.Range("A5").Value = "Kill"
$ws.Range("A6").Value = "Dustin"

# Fill column B for the remaining rows (reuses the "Adam" shared string)
$ws.Range("B4").Value = "Adam"
$ws.Range("B5").Value = "Adam"
$ws.Range("B6").Value = "Adam"

# Fill column C for the remaining rows (reuses the "ydg@hotmail.com" shared string)
$ws.Range("C4").Value = "ydg@hotmail.com"
$ws.Range("C5").Value = "ydg@hotmail.com"
$ws.Range("C6").Value = "ydg@hotmail.com"

# New course values first (introduces new shared strings)
$ws.Range("D3").Value = "Selenium Automation"
$ws.Range("D5").Value = "Mobile Automation"

# Reused course values
$ws.Range("D4").Value = "Agile & Manual Testing"
$ws.Range("D6").Value = "Agile & Manual Testing"

# Column E, "Others" for every new row
$ws.Range("E3").Value = "Others"
$ws.Range("E4").Value = "Others"
$ws.Range("E5").Value = "Others"
$ws.Range("E6").Value = "Others"

# Hyperlink the email cells, same as the existing C2 mailto link
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:ydg@hotmail.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:ydg@hotmail.com")
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:ydg@hotmail.com")
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:ydg@hotmail.com")

# Restore the selected cell to E4 as in the final saved workbook
$ws.Range("E4").Select()
